$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as the new row 14, pushing
# all subsequent data rows (old 14-123) down by one (new 15-124).
$ws.Rows("14:14").Insert()

# The inserted row inherits the static/categorical fields (Mercado ID,
# Mercado, Region, Codreg, Categoria ID, Categoria, Variedad, Calidad,
# Origen, Precio $/Kg, Clasificacion) from the row that is now directly
# below it (row 15, which holds what used to be row 14's data).
$ws.Range("A15:R15").Copy($ws.Range("A14:R14"))

# Now overwrite the cells that differ for this new observation: Fecha,
# Volumen, Precio minimo/maximo/promedio, Unidad de comercializacion and
# Kg o Unidades.
$ws.Range("D14").Value = 44685
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("N14").Value = "$/caja 50 unidades"
$ws.Range("Q14").Value = 50
